$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.274.27'
$ws.Range('E2').Value = '  -1.60%  '
$ws.Range('D3').Value = '3.433.20'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '231.72'
$ws.Range('E5').Value = '  -3.06%  '
$ws.Range('D6').Value = '621.11'
$ws.Range('E6').Value = '  -3.57%  '
$ws.Range('D7').Value = '1.39'
$ws.Range('E7').Value = '  -3.89%  '
$ws.Range('D8').Value = '0.394'
$ws.Range('E8').Value = '  -3.19%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '0.966'
$ws.Range('E10').Value = '  -1.45%  '
$ws.Range('D11').Value = '3.431.19'
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('D12').Value = '43.01'
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').Value = '6.23'
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('D15').Value = '93.117.81'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').Value = '4.073.12'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('E18').Value = '  -2.40%  '
$ws.Range('D19').Value = '3.433.44'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').Value = '18.23'
$ws.Range('E20').Value = '  +3.60%  '
$ws.Range('D21').Value = '11.76'
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('D22').Value = '500.46'
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D23').Value = '3.35'
$ws.Range('E23').Value = '  +2.67%  '
$ws.Range('D24').Value = '0.443'
$ws.Range('E24').Value = '  -10.31%  '
$ws.Range('D25').Value = '6.67'
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('E26').Value = '  -4.48%  '
$ws.Range('D27').Value = '92.10'
$ws.Range('E27').Value = '  -2.56%  '
$ws.Range('D28').Value = '12.00'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = '3.607.82'
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '11.46'
$ws.Range('E30').Value = '  -2.56%  '
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('D32').Value = '2.74'
$ws.Range('E32').Value = '  -0.90%  '
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('D34').Value = '1.01'
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('E35').Value = '  -4.01%  '
$ws.Range('D36').Value = '30.18'
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('D37').Value = '0.546'
$ws.Range('E37').Value = '  -1.68%  '
$ws.Range('D38').Value = '550.61'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('E39').Value = '  -3.50%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '1.41'
$ws.Range('E40').Value = '  -3.62%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').Value = '0.924'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('D43').Value = '0.150'
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').Value = '1.73'
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = '23.69'
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('B46').Value = 'MantraDAO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D46').Value = '3.71'
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('D47').Value = '5.53'
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('E48').Value = '  -0.63%  '
$ws.Range('E49').Value = '  -2.82%  '
$ws.Range('E50').Value = '  -5.37%  '
$ws.Range('E51').Value = '  -0.53%  '

Write-Output "Applied 96 cell updates"
